# Fixed query issue for C3DC phs002599
# The "Treatment Agent" column in the TreatmentTab query wrapped the
# REPLACE(...) call in a redundant CONCAT(...) call. Remove the CONCAT
# wrapper so the cell reads REPLACE(trt.treatment_agent, ';', ', ') directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B5")
$query = $cell.Value2

$oldFragment = "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))"
$newFragment = "REPLACE(trt.treatment_agent, ';', ', ')"

$updatedQuery = $query.Replace($oldFragment, $newFragment)

$cell.Value = $updatedQuery

# Re-apply the cell's font (same visible size/name it already had) so the
# edited cell gets its own formatting record, matching how Excel tracks a
# freshly retyped/pasted cell.
$cell.Font.Bold = $false
$cell.Font.Size = 12
$cell.Font.Name = "Calibri"

# Leave the cursor/selection where the edited cell is, matching the
# post-edit view position.
$ws.Range("C5").Select()
